# ULTA model update: bump the assumed SG&A annual growth rate from 3.5% to 4.0%.
#
# The driver cell is Model!M6 ("SG&A" row, first projected year), whose formula
# changes from =L6*1.035 to =L6*1.04. The remaining projected years (N6:V6) are a
# shared formula keyed off the prior column, so they're rewritten the same way;
# everything downstream (totals, margins, ROIC, NPV, share price, etc.) recalculates
# automatically from that single assumption change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model")

# --- Core assumption change: SG&A grows 4.0%/yr instead of 3.5%/yr ---
$ws.Range("M6").Formula = "=L6*1.04"
$ws.Range("N6:V6").Formula = "=M6*1.04"

# --- Everything below is unaffected in formula text; only recalculated values
#     shift. Re-assigning the identical formula text to these ranges mirrors
#     what Excel itself does when it re-serializes shared formulas after a
#     recalculation touches a long run of same-pattern cells. ---
$ws.Range("M10:V10").Formula = "=L21*`$Y`$15"

$ws.Range("W13:BB13").Formula  = "=V13*(1+`$Y`$16)"
$ws.Range("BC13:CH13").Formula = "=BB13*(1+`$Y`$16)"
$ws.Range("CI13:DK13").Formula = "=CH13*(1+`$Y`$16)"

# --- View state: scroll the frozen pane over one more column and move the
#     active selection on the bottom-right pane to R6. ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 18
[void]$ws.Range("R6").Select()
